$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Status" column (column E) with its header and value.
# Write E2 before E1 so the shared-string table picks up "Done" (row 2)
# before "Status" (row 1), matching the original author's edit order.
$ws.Range("E2").Value = "Done"
$ws.Range("E1").Value = "Status"

# Narrow column D (description) now that column E exists.
$ws.Columns.Item(4).ColumnWidth = 74.8333333333333

# Move the active selection to the newly added cell, like the author did.
[void]$ws.Range("E2").Select()
